# Fixing call to build array predictions.
# Update the contest definition name, the pool name, and the
# start/end dates used to build the array of predictions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Contest name used by "Create Contest Definition"
$ws.Range("B3").Value = "Russia2018"

# Start date day (row 5): June 1 2018 -> June 14 2018
$ws.Range("D5").Value = 14

# End date month/day (row 6): July 1 2018 -> June 29 2018
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 29

# Pool name used by "Create Contest Pool"
$ws.Range("B15").Value = "Bitbrackets"

# Update the active selection to match the authored state
$ws.Range("C18").Select()
